$d = $word.ActiveDocument

# Locate "This is a Microsoft word document." via Find instead of relying on
# hard-coded character offsets, so the script is resilient to any other
# differences in the document.
$found = $d.Content.Duplicate
$ok = $found.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find the target sentence 'This is a Microsoft word document.'"
}

# $found is now collapsed to the matched text; its End is the insertion point
# right after the period (and before the paragraph mark).
$insertPoint = $found.End

# Word's Range.InsertAfter happily appends text right into the neighbouring
# run when the run formatting is identical, merging what should be separate
# runs into a single <w:r>. To reproduce the commit's three distinct runs
# (" (", "Changed main", ")") we bracket every insertion with a momentary
# bookmark: the bookmark forces a run boundary at the insertion point, and
# we delete the bookmark again right after, leaving no trace of it behind.
function Insert-RunAt($doc, $position, $text, $bookmarkName) {
    $doc.Bookmarks.Add($bookmarkName, $doc.Range($position, $position)) | Out-Null
    $r = $doc.Range($position, $position)
    $r.InsertAfter($text)
    $doc.Bookmarks($bookmarkName).Delete()
    return $r.End
}

$pos = $insertPoint
$pos = Insert-RunAt $d $pos " (" "zzTmpSplitMark1"
$pos = Insert-RunAt $d $pos "Changed main" "zzTmpSplitMark2"
$pos = Insert-RunAt $d $pos ")" "zzTmpSplitMark3"

Write-Output ("Paragraph now reads: " + $d.Paragraphs(1).Range.Text)
